# Update stream_id values (column D) for rows 2-6 on Sheet1, and leave the
# selection on the edited range, matching the author's interactive edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2:D6").Value = 10194

$ws.Range("D2:D6").Select()
